$d = $word.ActiveDocument

# --- 1) Paragraph: "The explanatory variable ..." ---
$p6 = $d.Paragraphs.Item(6)
$rng6 = $p6.Range
$rng6.Find.Execute("The explanatory variable are", $true, $false, $false, $false, $false, $true, 1, $false, "The explanatory variable is a group dummy. Those dummies are", 2) | Out-Null

$rng6b = $p6.Range
$rng6b.Find.Execute("women who have a vegetarian diet and women who eat meat.", $true, $false, $false, $false, $false, $true, 1, $false, "women with a vegetarian-, a pescetarian diet, women who eat meat sometimes and women who eat meat regularly.", 2) | Out-Null

# --- 2) Econometric model formula paragraph ---
$p14 = $d.Paragraphs.Item(14)
$om = $p14.Range.OMaths.Item(1)
$om.Range.Text = ""

$omathXml = @"
<m:oMathPara><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/><w:lang w:val="en-GB"/></w:rPr><m:t>f</m:t></m:r><m:d><m:dPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/><w:lang w:val="en-GB"/></w:rPr></m:ctrlPr></m:dPr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/><w:lang w:val="en-GB"/></w:rPr><m:t>x</m:t></m:r></m:e></m:d><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/><w:lang w:val="en-GB"/></w:rPr><m:t>=</m:t></m:r><m:sSub><m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/><w:lang w:val="en-GB"/></w:rPr></m:ctrlPr></m:sSubPr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/><w:lang w:val="en-GB"/></w:rPr><m:t>β</m:t></m:r></m:e><m:sub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/><w:lang w:val="en-GB"/></w:rPr><m:t>0</m:t></m:r></m:sub></m:sSub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/><w:lang w:val="en-GB"/></w:rPr><m:t>+</m:t></m:r><m:sSub><m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/><w:i/><w:lang w:val="en-GB"/></w:rPr></m:ctrlPr></m:sSubPr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/><w:lang w:val="en-GB"/></w:rPr><m:t>β</m:t></m:r></m:e><m:sub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/><w:lang w:val="en-GB"/></w:rPr><m:t>1</m:t></m:r></m:sub></m:sSub><w:proofErr w:type="spellStart"/><m:r><m:rPr><m:nor/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/><w:lang w:val="en-GB"/></w:rPr><m:t>occ.meat</m:t></m:r><w:proofErr w:type="spellEnd"/><m:r><m:rPr><m:nor/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/><w:lang w:val="en-GB"/></w:rPr><m:t xml:space="preserve"> + </m:t></m:r><m:sSub><m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/><w:i/><w:lang w:val="en-GB"/></w:rPr></m:ctrlPr></m:sSubPr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/><w:lang w:val="en-GB"/></w:rPr><m:t>β</m:t></m:r></m:e><m:sub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/><w:lang w:val="en-GB"/></w:rPr><m:t>2</m:t></m:r></m:sub></m:sSub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/><w:lang w:val="en-GB"/></w:rPr><m:t xml:space="preserve">pescetarian+ </m:t></m:r><m:sSub><m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/><w:i/><w:lang w:val="en-GB"/></w:rPr></m:ctrlPr></m:sSubPr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/><w:lang w:val="en-GB"/></w:rPr><m:t>β</m:t></m:r></m:e><m:sub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/><w:lang w:val="en-GB"/></w:rPr><m:t>3</m:t></m:r></m:sub></m:sSub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/><w:lang w:val="en-GB"/></w:rPr><m:t>vegetarian</m:t></m:r><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/><w:lang w:val="en-GB"/></w:rPr><m:t>+</m:t></m:r><m:sSub><m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/><w:i/><w:lang w:val="en-GB"/></w:rPr></m:ctrlPr></m:sSubPr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/><w:lang w:val="en-GB"/></w:rPr><m:t>β</m:t></m:r></m:e><m:sub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/><w:lang w:val="en-GB"/></w:rPr><m:t>4</m:t></m:r></m:sub></m:sSub><m:r><m:rPr><m:nor/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/><w:lang w:val="en-GB"/></w:rPr><m:t>age</m:t></m:r><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/><w:lang w:val="en-GB"/></w:rPr><m:t>+ε</m:t></m:r></m:oMath></m:oMathPara>

"@

$pkgXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math">' + '<w:body><w:p>' + $omathXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$om2 = $p14.Range.OMaths.Item(1)
$om2.Range.InsertXML($pkgXml) | Out-Null

# --- 3) Paragraph with just "The" -> append headline sentence ---
$p16 = $d.Paragraphs.Item(16)
$p16.Range.InsertAfter(" headline of the article claims a clear relationship of breaking hips and vegetarian diet. The wording in the articles text is a bit more cautious. ")
